# Update handback datetime values for the "Generate Report for handback" commit.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-20 03:55:35"
$wsZhCn.Range("G5").Value = "2016-01-20 03:56:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-20 03:55:45"
$wsDeDe.Range("G5").Value = "2016-01-20 03:56:35"
